$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44253
$ws.Cells.Item(2, 10).Value = 95
$ws.Cells.Item(2, 11).Value = 9500
$ws.Cells.Item(2, 12).Value = 10000
$ws.Cells.Item(2, 13).Value = 9658
$ws.Cells.Item(2, 15).Value = "Región del Maule"
$ws.Cells.Item(2, 16).Value = 161

# Row 3
$ws.Cells.Item(3, 4).Value = 44266
$ws.Cells.Item(3, 10).Value = 60
$ws.Cells.Item(3, 11).Value = 9000
$ws.Cells.Item(3, 12).Value = 9500
$ws.Cells.Item(3, 13).Value = 9208
$ws.Cells.Item(3, 15).Value = "Región del Maule"
$ws.Cells.Item(3, 16).Value = 153

# Row 4
$ws.Cells.Item(4, 4).Value = 44202
$ws.Cells.Item(4, 10).Value = 50
$ws.Cells.Item(4, 11).Value = 8000
$ws.Cells.Item(4, 12).Value = 9000
$ws.Cells.Item(4, 13).Value = 8400
$ws.Cells.Item(4, 15).Value = "Región del Maule"
$ws.Cells.Item(4, 16).Value = 140

# Row 5
$ws.Cells.Item(5, 4).Value = 44694
$ws.Cells.Item(5, 10).Value = 80
$ws.Cells.Item(5, 11).Value = 10000
$ws.Cells.Item(5, 12).Value = 10000
$ws.Cells.Item(5, 13).Value = 10000
$ws.Cells.Item(5, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(5, 16).Value = 167

# Row 6
$ws.Cells.Item(6, 4).Value = 44627
$ws.Cells.Item(6, 10).Value = 60
$ws.Cells.Item(6, 11).Value = 9000
$ws.Cells.Item(6, 12).Value = 9500
$ws.Cells.Item(6, 13).Value = 9250
$ws.Cells.Item(6, 15).Value = "Región Metropolitana"
$ws.Cells.Item(6, 16).Value = 154

# Row 7
$ws.Cells.Item(7, 4).Value = 44216
$ws.Cells.Item(7, 10).Value = 55
$ws.Cells.Item(7, 11).Value = 9500
$ws.Cells.Item(7, 12).Value = 10000
$ws.Cells.Item(7, 13).Value = 9773
$ws.Cells.Item(7, 15).Value = "Región del Maule"
$ws.Cells.Item(7, 16).Value = 163

# Row 8
$ws.Cells.Item(8, 4).Value = 44600
$ws.Cells.Item(8, 10).Value = 60
$ws.Cells.Item(8, 11).Value = 12000
$ws.Cells.Item(8, 12).Value = 13000
$ws.Cells.Item(8, 13).Value = 12500
$ws.Cells.Item(8, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 16).Value = 208

# Row 9
$ws.Cells.Item(9, 4).Value = 44264
$ws.Cells.Item(9, 10).Value = 43
$ws.Cells.Item(9, 11).Value = 8500
$ws.Cells.Item(9, 12).Value = 9000
$ws.Cells.Item(9, 13).Value = 8709
$ws.Cells.Item(9, 15).Value = "Región del Maule"
$ws.Cells.Item(9, 16).Value = 145

# Row 10
$ws.Cells.Item(10, 4).Value = 44610
$ws.Cells.Item(10, 10).Value = 60
$ws.Cells.Item(10, 11).Value = 11000
$ws.Cells.Item(10, 12).Value = 12000
$ws.Cells.Item(10, 13).Value = 11500
$ws.Cells.Item(10, 15).Value = "Región Metropolitana"
$ws.Cells.Item(10, 16).Value = 192

# Row 11
$ws.Cells.Item(11, 4).Value = 44594
$ws.Cells.Item(11, 10).Value = 80
$ws.Cells.Item(11, 11).Value = 12000
$ws.Cells.Item(11, 12).Value = 13000
$ws.Cells.Item(11, 13).Value = 12500
$ws.Cells.Item(11, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 16).Value = 208

# Row 12
$ws.Cells.Item(12, 4).Value = 44687
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 9000
$ws.Cells.Item(12, 12).Value = 10000
$ws.Cells.Item(12, 13).Value = 9500
$ws.Cells.Item(12, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(12, 16).Value = 158

# Row 13
$ws.Cells.Item(13, 4).Value = 44224
$ws.Cells.Item(13, 10).Value = 80
$ws.Cells.Item(13, 11).Value = 8500
$ws.Cells.Item(13, 12).Value = 9000
$ws.Cells.Item(13, 13).Value = 8719
$ws.Cells.Item(13, 15).Value = "Región del Maule"
$ws.Cells.Item(13, 16).Value = 145

# Row 14
$ws.Cells.Item(14, 4).Value = 44671
$ws.Cells.Item(14, 10).Value = 160
$ws.Cells.Item(14, 11).Value = 6000
$ws.Cells.Item(14, 12).Value = 7000
$ws.Cells.Item(14, 13).Value = 6500
$ws.Cells.Item(14, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(14, 16).Value = 108

# Row 15
$ws.Cells.Item(15, 4).Value = 44271
$ws.Cells.Item(15, 10).Value = 55
$ws.Cells.Item(15, 11).Value = 9000
$ws.Cells.Item(15, 12).Value = 9500
$ws.Cells.Item(15, 13).Value = 9227
$ws.Cells.Item(15, 15).Value = "Región del Maule"
$ws.Cells.Item(15, 16).Value = 154

# Row 16
$ws.Cells.Item(16, 4).Value = 44624
$ws.Cells.Item(16, 10).Value = 60
$ws.Cells.Item(16, 11).Value = 10000
$ws.Cells.Item(16, 12).Value = 11000
$ws.Cells.Item(16, 13).Value = 10500
$ws.Cells.Item(16, 15).Value = "Región Metropolitana"
$ws.Cells.Item(16, 16).Value = 175

# Row 17
$ws.Cells.Item(17, 4).Value = 44259
$ws.Cells.Item(17, 10).Value = 70
$ws.Cells.Item(17, 11).Value = 9000
$ws.Cells.Item(17, 12).Value = 9500
$ws.Cells.Item(17, 13).Value = 9214
$ws.Cells.Item(17, 15).Value = "Región del Maule"
$ws.Cells.Item(17, 16).Value = 154

# Row 18
$ws.Cells.Item(18, 4).Value = 44608
$ws.Cells.Item(18, 10).Value = 100
$ws.Cells.Item(18, 11).Value = 12000
$ws.Cells.Item(18, 12).Value = 13000
$ws.Cells.Item(18, 13).Value = 12500
$ws.Cells.Item(18, 15).Value = "Región Metropolitana"
$ws.Cells.Item(18, 16).Value = 208

# Row 19
$ws.Cells.Item(19, 4).Value = 44630
$ws.Cells.Item(19, 10).Value = 60
$ws.Cells.Item(19, 11).Value = 9000
$ws.Cells.Item(19, 12).Value = 9500
$ws.Cells.Item(19, 13).Value = 9250
$ws.Cells.Item(19, 15).Value = "Región Metropolitana"
$ws.Cells.Item(19, 16).Value = 154

# Row 20
$ws.Cells.Item(20, 4).Value = 44204
$ws.Cells.Item(20, 10).Value = 45
$ws.Cells.Item(20, 11).Value = 9500
$ws.Cells.Item(20, 12).Value = 10000
$ws.Cells.Item(20, 13).Value = 9722
$ws.Cells.Item(20, 15).Value = "Región del Maule"
$ws.Cells.Item(20, 16).Value = 162

# Row 21
$ws.Cells.Item(21, 4).Value = 44160
$ws.Cells.Item(21, 10).Value = 90
$ws.Cells.Item(21, 11).Value = 7500
$ws.Cells.Item(21, 12).Value = 8000
$ws.Cells.Item(21, 13).Value = 7667
$ws.Cells.Item(21, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(21, 16).Value = 128

# Row 22
$ws.Cells.Item(22, 4).Value = 44698
$ws.Cells.Item(22, 10).Value = 60
$ws.Cells.Item(22, 11).Value = 10000
$ws.Cells.Item(22, 12).Value = 10000
$ws.Cells.Item(22, 13).Value = 10000
$ws.Cells.Item(22, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(22, 16).Value = 167

# Row 23
$ws.Cells.Item(23, 4).Value = 44692
$ws.Cells.Item(23, 10).Value = 60
$ws.Cells.Item(23, 11).Value = 10000
$ws.Cells.Item(23, 12).Value = 10000
$ws.Cells.Item(23, 13).Value = 10000
$ws.Cells.Item(23, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(23, 16).Value = 167

# Row 24
$ws.Cells.Item(24, 4).Value = 44617
$ws.Cells.Item(24, 10).Value = 100
$ws.Cells.Item(24, 11).Value = 10000
$ws.Cells.Item(24, 12).Value = 11000
$ws.Cells.Item(24, 13).Value = 10500
$ws.Cells.Item(24, 15).Value = "Región Metropolitana"
$ws.Cells.Item(24, 16).Value = 175

# Row 25
$ws.Cells.Item(25, 4).Value = 44615
$ws.Cells.Item(25, 10).Value = 100
$ws.Cells.Item(25, 11).Value = 11000
$ws.Cells.Item(25, 12).Value = 12000
$ws.Cells.Item(25, 13).Value = 11500
$ws.Cells.Item(25, 15).Value = "Región Metropolitana"
$ws.Cells.Item(25, 16).Value = 192

# Row 26
$ws.Cells.Item(26, 4).Value = 44208
$ws.Cells.Item(26, 10).Value = 100
$ws.Cells.Item(26, 11).Value = 7000
$ws.Cells.Item(26, 12).Value = 8000
$ws.Cells.Item(26, 13).Value = 7350
$ws.Cells.Item(26, 15).Value = "Región del Maule"
$ws.Cells.Item(26, 16).Value = 122

# Row 27
$ws.Cells.Item(27, 4).Value = 44210
$ws.Cells.Item(27, 10).Value = 60
$ws.Cells.Item(27, 11).Value = 8000
$ws.Cells.Item(27, 12).Value = 9000
$ws.Cells.Item(27, 13).Value = 8417
$ws.Cells.Item(27, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(27, 16).Value = 140

# Row 28
$ws.Cells.Item(28, 4).Value = 44218
$ws.Cells.Item(28, 10).Value = 65
$ws.Cells.Item(28, 11).Value = 9000
$ws.Cells.Item(28, 12).Value = 10000
$ws.Cells.Item(28, 13).Value = 9615
$ws.Cells.Item(28, 15).Value = "Región del Maule"
$ws.Cells.Item(28, 16).Value = 160

# Row 29
$ws.Cells.Item(29, 4).Value = 44162
$ws.Cells.Item(29, 10).Value = 43
$ws.Cells.Item(29, 11).Value = 8000
$ws.Cells.Item(29, 12).Value = 8500
$ws.Cells.Item(29, 13).Value = 8209
$ws.Cells.Item(29, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(29, 16).Value = 137

# Row 30
$ws.Cells.Item(30, 4).Value = 44699
$ws.Cells.Item(30, 10).Value = 60
$ws.Cells.Item(30, 11).Value = 10000
$ws.Cells.Item(30, 12).Value = 10000
$ws.Cells.Item(30, 13).Value = 10000
$ws.Cells.Item(30, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(30, 16).Value = 167

# Row 31
$ws.Cells.Item(31, 4).Value = 44596
$ws.Cells.Item(31, 10).Value = 100
$ws.Cells.Item(31, 11).Value = 12000
$ws.Cells.Item(31, 12).Value = 13000
$ws.Cells.Item(31, 13).Value = 12500
$ws.Cells.Item(31, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(31, 16).Value = 208

# Row 32
$ws.Cells.Item(32, 4).Value = 44159
$ws.Cells.Item(32, 10).Value = 35
$ws.Cells.Item(32, 11).Value = 7500
$ws.Cells.Item(32, 12).Value = 8000
$ws.Cells.Item(32, 13).Value = 7714
$ws.Cells.Item(32, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(32, 16).Value = 129

# Row 33
$ws.Cells.Item(33, 4).Value = 44690
$ws.Cells.Item(33, 10).Value = 60
$ws.Cells.Item(33, 11).Value = 10000
$ws.Cells.Item(33, 12).Value = 10000
$ws.Cells.Item(33, 13).Value = 10000
$ws.Cells.Item(33, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(33, 16).Value = 167
